$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-02-14"

# Update the "2022 (through 02-13)" header label to "02-14"
$ws.Range("I1").Value = "2022 (through 02-14)"

# Update the February (row 3) 2022-to-date count
$ws.Range("I3").Value = 65

# Update the Total row (row 14) 2022-to-date count
$ws.Range("I14").Value = 226
